$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are numeric-looking text (e.g. "0", "1057.59...", "-3.48e-05")
# must keep their original text storage type, so force Text number format before
# assigning the value (otherwise Excel auto-converts the string to a real number).
$numericLookingCells = @(
    "I11", "K11", "I13", "K13", "J14", "L14", "I16", "K16", "J23", "L23", "I24", "K24", "I28", "K28", "I33", "K33", "I37", "K37", "I39", "K39", "J44", "L44", "J48", "L48"
)
foreach ($ref in $numericLookingCells) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply the new cell values from the "3 Pt Linearity Test" results update.
$updates = @{
    'C2' = 'PASS'
    'C3' = 'PASS'
    'C4' = 'PASS'
    'G4' = 'LIN'
    'C5' = 'PASS'
    'G5' = 'LIN'
    'C6' = 'PASS'
    'G6' = 'LIN'
    'C7' = 'PASS'
    'G7' = 'LIN'
    'C8' = 'PASS'
    'G8' = 'LIN'
    'C9' = 'PASS'
    'G9' = 'LIN'
    'C10' = 'PASS'
    'C11' = 'PASS'
    'I11' = '0'
    'K11' = '0'
    'C12' = 'PASS'
    'C13' = 'PASS'
    'I13' = '0'
    'K13' = '0'
    'C14' = 'PASS'
    'H14' = 'LIN'
    'J14' = '0.000222388359832818'
    'L14' = '1057.59008402095'
    'C15' = 'PASS'
    'G15' = 'LIN'
    'C16' = 'PASS'
    'I16' = '0'
    'K16' = '0'
    'C17' = 'PASS'
    'G17' = 'LIN'
    'C18' = 'PASS'
    'C19' = 'PASS'
    'G19' = 'LIN'
    'C20' = 'PASS'
    'C21' = 'PASS'
    'G21' = 'LIN'
    'C22' = 'PASS'
    'G22' = 'LIN'
    'C23' = 'PASS'
    'H23' = 'LIN'
    'J23' = '0.000231374756624682'
    'L23' = '2038.39309588292'
    'C24' = 'PASS'
    'I24' = '0'
    'K24' = '0'
    'C25' = 'PASS'
    'C26' = 'PASS'
    'G26' = 'LIN'
    'C27' = 'PASS'
    'C28' = 'PASS'
    'I28' = '0'
    'K28' = '0'
    'C29' = 'PASS'
    'C30' = 'PASS'
    'C31' = 'PASS'
    'G31' = 'LIN'
    'C32' = 'PASS'
    'C33' = 'PASS'
    'I33' = '0'
    'K33' = '0'
    'C34' = 'PASS'
    'G34' = 'LIN'
    'C35' = 'PASS'
    'G35' = 'LIN'
    'C36' = 'PASS'
    'G36' = 'LIN'
    'C37' = 'PASS'
    'I37' = '0'
    'K37' = '0'
    'C38' = 'PASS'
    'G38' = 'LIN'
    'C39' = 'PASS'
    'I39' = '0'
    'K39' = '0'
    'C40' = 'PASS'
    'G40' = 'LIN'
    'C41' = 'PASS'
    'G41' = 'LIN'
    'C42' = 'PASS'
    'G42' = 'LIN'
    'C43' = 'PASS'
    'G43' = 'LIN'
    'C44' = 'PASS'
    'G44' = 'LIN'
    'H44' = 'LIN'
    'J44' = '-3.48379580108469e-05'
    'L44' = '-370.8319472873'
    'C45' = 'PASS'
    'G45' = 'LIN'
    'C46' = 'PASS'
    'C47' = 'PASS'
    'C48' = 'PASS'
    'H48' = 'LIN'
    'J48' = '3.98039650706417e-05'
    'L48' = '429.424281085316'
    'C49' = 'PASS'
    'G49' = 'LIN'
}
foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
